# Generate Report for Handback
# Updates the handoff/handback timestamp strings on the report sheets to
# reflect a newer report-generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-09-06 11:27:12"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-06 11:26:59"

# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-09-06 11:27:32"

# de-de!H2 - "Correspond Handoff Datetime" (shares old value with Overview!G2)
$wsDeDe.Range("H2").Value = "2016-09-06 11:27:12"

# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-09-06 11:27:40"
